# Weekly refresh of the "Fruta, Agricola del Norte S.A. de Arica - Kiwi" price sheet.
# Updates Fecha/Calidad/Volumen/Precio minimo/Precio maximo/Precio promedio ponderado/
# Unidad de comercializacion/Precio $ por Kg/Kg por unidad for each data row (2-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns per row: Row, D(Fecha), L(Calidad), M(Volumen), N(Precio minimo), O(Precio maximo),
#                  P(Precio promedio ponderado), Q(Unidad de comercializacion), S(Precio $/Kg), T(Kg/unidad)
$data = @()
$data += ,@(2, 44656, 'Primera', 270, 19000, 20000, 19500, '$/bandeja 18 kilos', 1083, 18)
$data += ,@(3, 44819, 'Primera', 300, 17000, 18000, 17500, '$/bandeja 10 kilos', 1750, 10)
$data += ,@(4, 44307, 'Primera', 250, 19000, 20000, 19500, '$/bandeja 18 kilos', 1083, 18)
$data += ,@(5, 44991, 'Primera', 250, 24000, 25000, 24500, '$/bandeja 18 kilos', 1361, 18)
$data += ,@(6, 44491, 'Primera', 300, 14000, 15000, 14500, '$/bandeja 10 kilos', 1450, 10)
$data += ,@(7, 44784, 'Primera', 300, 19000, 20000, 19500, '$/bandeja 18 kilos', 1083, 18)
$data += ,@(8, 44614, 'Primera', 250, 20000, 21000, 20500, '$/bandeja 18 kilos', 1139, 18)
$data += ,@(9, 44323, 'Primera', 270, 21000, 22000, 21500, '$/bandeja 18 kilos', 1194, 18)
$data += ,@(10, 44291, 'Primera', 200, 17000, 18000, 17500, '$/bandeja 18 kilos', 972, 18)
$data += ,@(11, 44629, 'Segunda', 300, 17000, 18000, 17500, '$/bandeja 18 kilos', 972, 18)
$data += ,@(12, 44602, 'Primera', 270, 20000, 21000, 20500, '$/bandeja 18 kilos', 1139, 18)
$data += ,@(13, 44487, 'Primera', 300, 14000, 15000, 14500, '$/bandeja 10 kilos', 1450, 10)
$data += ,@(14, 44706, 'Primera', 400, 9000, 10000, 9500, '$/bandeja 10 kilos', 950, 10)
$data += ,@(15, 44263, 'Primera', 250, 21000, 22000, 21500, '$/caja 18 kilos', 1194, 18)
$data += ,@(16, 44616, 'Segunda', 300, 16000, 17000, 16500, '$/caja 18 kilos granel', 917, 18)
$data += ,@(17, 44418, 'Primera', 240, 10000, 11000, 10500, '$/bandeja 10 kilos', 1050, 10)
$data += ,@(18, 44489, 'Primera', 300, 26000, 27000, 26500, '$/bandeja 18 kilos', 1472, 18)
$data += ,@(19, 44789, 'Segunda', 250, 19000, 20000, 19500, '$/bandeja 18 kilos', 1083, 18)
$data += ,@(20, 44673, 'Especial', 400, 14000, 15000, 14500, '$/bandeja 10 kilos', 1450, 10)
$data += ,@(21, 45002, 'Segunda', 300, 24000, 25000, 24500, '$/bandeja 18 kilos', 1361, 18)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D Fecha
    $ws.Cells.Item($r, 12).Value = $row[2]   # L Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]   # M Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[5]   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[6]   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[7]   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $row[8]   # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[9]   # T Kg / unidad
}
